# "finished mansion (for now)" - log a new entry for the mansion basement work.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 34: date 2025-01-10 (serial 45667), task description, hours worked.
$ws.Range("A34").Value = 45667
# Match the date formatting already used by the column (style applied to A2:A33).
$ws.Range("A34").NumberFormat = $ws.Range("A33").NumberFormat
$ws.Range("B34").Value = "mansion basement"
$ws.Range("C34").Value = 5

# Leave the selection on the next empty row, as recorded in the saved view.
$ws.Range("E35").Select()
